# "working on league mode" - add three new template rows (105, 106, 107)
# to the Templates sheet, describing new league-mode related UI steps.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Templates")

# Row 102 -> Template 105
$ws.Range("A102").Value = "105"
$ws.Range("B102").Value = 368
$ws.Range("C102").Value = 736
$ws.Range("D102").Value = 504
$ws.Range("E102").Value = 856
$ws.Range("F102").Value = "105"
$ws.Range("G102").Value = "half time stamina recovery +30% -> league mode"

# Row 103 -> Template 106
$ws.Range("A103").Value = "106"
$ws.Range("B103").Value = 718
$ws.Range("C103").Value = 736
$ws.Range("D103").Value = 854
$ws.Range("E103").Value = 856
$ws.Range("F103").Value = "106"
$ws.Range("G103").Value = "stamina consumed -20% -> league mode"

# Row 104 -> Template 107
$ws.Range("A104").Value = "107"
$ws.Range("B104").Value = 1357
$ws.Range("C104").Value = 718
$ws.Range("D104").Value = 1708
$ws.Range("E104").Value = 828
$ws.Range("F104").Value = "107"
$ws.Range("G104").Value = "kick off -> league mode"

# Match the author's final selection after entering the data
$ws.Range("G104").Select()
